$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-Rows($ws, $r1, $r2, $cols) {
    foreach ($col in $cols) {
        $addr1 = "$col$r1"
        $addr2 = "$col$r2"
        $v1 = $ws.Range($addr1).Value2
        $v2 = $ws.Range($addr2).Value2
        $ws.Range($addr1).Value = $v2
        $ws.Range($addr2).Value = $v1
    }
}

$cols = @("B", "C", "D", "E", "F", "G")

Swap-Rows $ws 76 77 $cols
Swap-Rows $ws 82 83 $cols
Swap-Rows $ws 86 87 $cols
Swap-Rows $ws 100 101 $cols
Swap-Rows $ws 189 190 $cols
Swap-Rows $ws 234 235 $cols
Swap-Rows $ws 366 367 $cols
Swap-Rows $ws 370 371 $cols
Swap-Rows $ws 409 410 $cols
Swap-Rows $ws 415 416 $cols
Swap-Rows $ws 417 418 $cols
Swap-Rows $ws 427 428 $cols
Swap-Rows $ws 497 498 $cols
Swap-Rows $ws 531 532 $cols
Swap-Rows $ws 535 536 $cols
Swap-Rows $ws 541 542 $cols
Swap-Rows $ws 748 749 $cols
Swap-Rows $ws 776 777 $cols
Swap-Rows $ws 778 779 $cols
Swap-Rows $ws 782 783 $cols
Swap-Rows $ws 784 785 $cols
Swap-Rows $ws 788 789 $cols
Swap-Rows $ws 872 873 $cols
Swap-Rows $ws 884 885 $cols
Swap-Rows $ws 887 888 $cols
Swap-Rows $ws 896 897 $cols
